$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the D55:D58 "idPagamento" values that were previously blank ---
$ws.Range("D55").NumberFormat = "@"
$ws.Range("D55").Value = "77136751137"

$ws.Range("D56").NumberFormat = "@"
$ws.Range("D56").Value = "77136970483"

$ws.Range("D57").NumberFormat = "@"
$ws.Range("D57").Value = "77381441200"

$ws.Range("D58").NumberFormat = "@"
$ws.Range("D58").Value = "77381528772"

# --- Append four new rows (61-64) with the corrected per-user number selections ---
function Add-SelectionRow {
    param($Row, $Name, $Id, $Phone, $Numbers)

    $ws.Range("A$Row").Value = $Name
    $ws.Range("B$Row").Value = $Id

    $ws.Range("C$Row").NumberFormat = "@"
    $ws.Range("C$Row").Value = $Phone

    # "idPagamento" (D) is left blank for these new rows, same as the source edit.

    $cols = @("E","F","G","H","I","J","K","L","M","N")
    for ($i = 0; $i -lt $Numbers.Length; $i++) {
        $ws.Range("$($cols[$i])$Row").Value = $Numbers[$i]
    }

    $ws.Range("O$Row").Value = "Não"
}

Add-SelectionRow 61 "Isabelly Silva Quintans" 7117522682 "11966548087" @(1,2,3,4,5,6,7,8,9,10)
Add-SelectionRow 62 "Vitor Ito" 1578424633 "11987876543" @(1,2,3,4,5,6,7,8,9,10)
Add-SelectionRow 63 "Vitor Ito" 1578424633 "11987541236" @(21,26,27,41,43,45,46,49,51,56)
Add-SelectionRow 64 "Isabelly Silva Quintans" 7117522682 "11987541236" @(21,26,27,41,43,45,46,49,51,56)

# --- Keep the "number stored as text" warning suppressed over the (now larger) used range ---
$usedRange = $ws.Range("A1:O64")
$numberAsTextError = $usedRange.Errors().Item(3)
$numberAsTextError.Ignore = $true
